$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B5: separator-config count changed from 3 to 0
$ws.Range("B5").Value = 0

# C14: clear cell value (field now omitted / blank)
$ws.Range("C14").ClearContents()
